$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1): columns C, D, E get relabeled
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Update data rows 2..12: column C becomes the family text, D stays the same
# family text, and E becomes numeric 1 (instead of text family)
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 3).Value = "f__Akkermansiaceae"
    $ws.Cells.Item($r, 5).Value = 1
}
